# Update dashboards - 2026-02-05
# Applies weekly refresh to the "Aguilar Prototype" sheet:
#  - advances the "Latest Date" cells for several series
#  - shifts the Present/Lag1/Lag2/Lag3/Lag4 (Q:U) values one column to the
#    right (Lag4 drops off, a new Present value is inserted)
#  - clears the "new data" yellow highlight (style 49 -> 48) on date cells
#    that are no longer the most-recently-updated row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Shift-Lags {
    param(
        [string]$RowNum,
        [double]$NewPresent
    )
    # Values currently in Q:T become R:U (U / old Lag4 is discarded),
    # and the new "present" value is written into Q.
    $q = $ws.Range("Q$RowNum").Value2
    $r = $ws.Range("R$RowNum").Value2
    $s = $ws.Range("S$RowNum").Value2
    $t = $ws.Range("T$RowNum").Value2

    $ws.Range("U$RowNum").Value = $t
    $ws.Range("T$RowNum").Value = $s
    $ws.Range("S$RowNum").Value = $r
    $ws.Range("R$RowNum").Value = $q
    $ws.Range("Q$RowNum").Value = $NewPresent
}

# ---- Row 13 : UI Initial Claims (ICSA) ----
$ws.Range("N13").Value = 46048
Shift-Lags "13" 231000

# ---- Row 14 : UI Continuing Claims (CCSA) ----
$ws.Range("N14").Value = 46041
Shift-Lags "14" 1844000

# ---- Row 29 : 5yr, 5yr Forward (T5YIFR) ----
$ws.Range("N29").Value = 46057
Shift-Lags "29" 2.19

# ---- Row 30 : 10yr TIPS (T10YIE) ----
$ws.Range("N30").Value = 46057
Shift-Lags "30" 2.35

# ---- Rows 46-51 : clear the "new data" yellow highlight on date cells ----
# These date cells keep their value but lose the yellow fill (style 49 -> 48).
# A cell that already carries style 48 (no fill) is copied and pasted as
# "Formats only" so the existing style index is reused rather than a new one
# being created.
$styleSource = $ws.Range("C3")
$styleSource.Copy() | Out-Null
foreach ($addr in @("C46","C47","C48","C49","C50","C51","N51")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# ---- Row 47 : FFR (DFF) ----
$ws.Range("N47").Value = 46056

# ---- Row 48 : 2y UST (DGS2) ----
$ws.Range("N48").Value = 46056
Shift-Lags "48" 3.57

# ---- Row 49 : 5y UST (DGS5) ----
$ws.Range("N49").Value = 46056
Shift-Lags "49" 3.83

# ---- Row 50 : 10y UST (DGS10) ----
$ws.Range("N50").Value = 46056
Shift-Lags "50" 4.28

# ---- Row 52 : BAA (DBAA) ----
$ws.Range("N52").Value = 46056
Shift-Lags "52" 5.91
